# Applies the "Optimizations returned, updated report" commit to the
# "PR3 Measurements" workbook (sheet "Blad1" / sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- New commentary notes (become new shared strings, H2/H3) ---
$ws.Range("H2").Value = " Results aren't 100% accurate but are decent indications"
$ws.Range("H3").Value = "Averaging out the results gives measurement errors, in DrawTanks for example"

# --- MAXP1 = 16 block ---
$ws.Range("B6").Value = 88104
$ws.Range("B7").Value = 95998558
$ws.Range("B8").Value = 95837864
$ws.Range("B9").Value = 162579

# --- MAXP1 = 32 block ---
$ws.Range("B14").Value = 307184
$ws.Range("B15").Value = 191395432
$ws.Range("B16").Value = 191121105
$ws.Range("B17").Value = 340570
$ws.Range("B18").Value = 8

# --- MAXP1 = 128 block ---
$ws.Range("B22").Value = 4666242
$ws.Range("B23").Value = 699634715
$ws.Range("B24").Value = 699763880
$ws.Range("B25").Value = 1066045

# --- MAXP1 = 256 block (previously-empty cells now filled in) ---
$ws.Range("B30").Value = 17914301
$ws.Range("B31").Value = 1384320501
$ws.Range("B32").Value = 138525268
$ws.Range("B33").Value = 1558078
$ws.Range("B34").Value = 1

# --- Selection / view state: drop the frozen topLeftCell scroll position
# and move the active selection from D33 to A26 ---
$ws.Range("A26").Select()
